$wb = $excel.ActiveWorkbook

# zh-cn sheet: update the "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" timestamps in row 2 (regenerated report)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 10:54:19"
$wsZhCn.Range("H2").Value = "2016-03-11 10:54:36"

# de-de sheet: same two timestamps in row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 10:54:22"
$wsDeDe.Range("H2").Value = "2016-03-11 10:54:41"
